$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.266.68"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.872.42"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "'241.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'0.3107"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'0.07701"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "'25.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "'0.08382"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "1.885.75"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'5.211"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'0.7104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "29.271.00"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000008272"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.12%  "
$ws.Range("D18").Value = "'5.961"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'242.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "2.128.95"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'7.818"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'0.1629"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("D26").Value = "'163.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").Value = "'9.010"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").Value = "'18.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "'1.503"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "'4.412"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'4.318"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.42%  "
$ws.Range("D32").Value = "'1.283"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.29%  "
$ws.Range("D33").Value = "'0.05242"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "'1.923"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.171"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7465"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "'2.682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.01855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "'2.718"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").Value = "1.154.92"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'73.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "'0.8856"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").Value = "'105.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").Value = "'0.9995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "2.024.84"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'1.801"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "'0.5188"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'9.373"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'0.4297"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
